$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right above the existing row 360, pushing all the
# current rows 360-459 down to 363-462 (matching the target diff, which
# shows a brand-new date block - 2021-11-08 / serial 44508 - inserted
# ahead of the existing "44333" block, with everything after it unchanged
# in relative order).
$ws.Rows.Item(360).Insert()
$ws.Rows.Item(360).Insert()
$ws.Rows.Item(360).Insert()

# Populate the 3 freshly inserted rows (360-362) with the new price
# records for Terminal La Palmera de La Serena - Frutilla, date 44508.
$newRows = @(
    @{ Row = 360; Calidad = "Especial"; Volumen = 400; Min = 12500; Max = 13000; Prom = 12750; PrecioKg = 1821 },
    @{ Row = 361; Calidad = "Primera";  Volumen = 400; Min = 10500; Max = 11000; Prom = 10750; PrecioKg = 1536 },
    @{ Row = 362; Calidad = "Segunda";  Volumen = 360; Min = 8500;  Max = 9000;  Prom = 8750;  PrecioKg = 1250 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44508
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "$/bandeja 7 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 7
}
